$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.297.01"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.028.32"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'227.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'55.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.32%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").Value = "'0.0793"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  -5.12%  "
$ws.Range("D12").Value = "2.327.29"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").Value = "'20.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").Value = "2.030.05"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "37.223.73"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'6.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("D20").Value = "'69.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "'224.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  -6.19%  "
$ws.Range("D26").Value = "'9.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.68%  "
$ws.Range("D27").Value = "'165.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("D28").Value = "'0.129"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("D32").Value = "'4.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("E35").Value = "  -5.60%  "
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -4.41%  "
$ws.Range("D39").Value = "'5.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D41").Value = "1.477.94"
$ws.Range("D42").Value = "'96.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("D43").Value = "'16.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").Value = "'2.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("E46").Value = "  -4.97%  "
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "'7.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "2.215.22"
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("E51").Value = "  -9.94%  "
